$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Validation expression for "date record was created" (row 7)
$dateRecordWasCreatedValidation = @"
(function() {
  if (dateRecordWasCreated !== null) {
    const today = new Date();
    const earliestDate = new Date("2015-01-01");
    const input = new Date(dateRecordWasCreated);
    const todayParsed = Date.parse(today);
    const inputParsed = Date.parse(input);
    const earliestParsed = Date.parse(earliestDate);
    if (inputParsed < earliestParsed) {
      return "The date the record was created cannot be earlier than 01 January 2015";
    }
    if (inputParsed > todayParsed) {
      return "The date the record was created cannot be greater than today";
    }
  }
})();
"@

# Validation expression for "date record was last modified" (row 8)
$dateRecordWasLastModifiedValidation = @"
(function() {
  if (dateRecordWasLastModified !== null) {
    const today = new Date();
    const earliestDate = new Date("2015-01-01");
    const input = new Date(dateRecordWasLastModified);
    const todayParsed = Date.parse(today);
    const inputParsed = Date.parse(input);
    const earliestParsed = Date.parse(earliestDate);
    if (inputParsed < earliestParsed) {
      return "The date the record was created cannot be earlier than 01 January 2015";
    }
    if (inputParsed > todayParsed) {
      return "The date the record was created cannot be greater than today";
    }
  }
})();
"@

# Validation expression for "date informed consent was given" (row 10)
$dateInformedConsentWasGivenValidation = @"
(function() {
  if (dateInformedConsentWasGiven !== null) {
    const today = new Date();
    const earliestDate = new Date("2015-01-01");
    const input = new Date(dateInformedConsentWasGiven);
    const todayParsed = Date.parse(today);
    const inputParsed = Date.parse(input);
    const earliestParsed = Date.parse(earliestDate);
    if (inputParsed < earliestParsed) {
      return "The date the record was created cannot be earlier than 01 January 2015";
    }
    if (inputParsed > todayParsed) {
      return "The date the record was created cannot be greater than today";
    }
  }
})();
"@

$ws.Range("J7").Value = $dateRecordWasCreatedValidation
$ws.Range("J8").Value = $dateRecordWasLastModifiedValidation
$ws.Range("J10").Value = $dateInformedConsentWasGivenValidation
